# Offline mode to get all variables - backlog sheet update.
# Add the new "Improve search and refresh" task at the top of the
# Medium/Low priority block, push the existing tasks in that block down by
# one row, and mark several already-delivered tasks as "Done".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

# --- Row 21: brand new task (was "to have a log file / ... error on search") ---
$ws.Range("C21").Value = "High"
$ws.Range("D21").Value = "To do"
$ws.Range("E21").Value = "Improve search and refresh"
$ws.Range("F21").Value = "to more fluidity and avoid undefined variables"
$ws.Range("G21").Value = 4
$ws.Range("H21").Value = $null

# --- Row 22: same task, now marked Done (will be filtered/hidden) ---
$ws.Range("C22").Value = "Medium"
$ws.Range("D22").Value = "Done"
$ws.Range("E22").Value = "to copy the path and name when we are in the path column"
$ws.Range("F22").Value = "to have directly the path and name"
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 0.0625

# --- Row 23: what used to be row 21's task ---
$ws.Range("C23").Value = "Medium"
$ws.Range("D23").Value = "To do"
$ws.Range("E23").Value = "to have a log file"
$ws.Range("F23").Value = "to see the diffrent application error"
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = $null

# --- Row 24: unchanged, left untouched on purpose ---

# --- Row 25: same task, now marked Done (will be filtered/hidden) ---
$ws.Range("C25").Value = "Medium"
$ws.Range("D25").Value = "Done"
$ws.Range("E25").Value = "the locked list are automatically saved"
$ws.Range("F25").Value = "to load automatically the locked list and the application shouldn't carsh"
$ws.Range("G25").Value = 3
$ws.Range("H25").Value = 1

# --- Row 26: same task, now marked Done (will be filtered/hidden) ---
$ws.Range("C26").Value = "Medium"
$ws.Range("D26").Value = "Done"
$ws.Range("E26").Value = "to implement the comment in XML file (rules coloring)"
$ws.Range("F26").Value = "to show the color's signification"
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 0.020833333333333332

# --- Row 27: what used to be row 23's task ---
$ws.Range("C27").Value = "Medium"
$ws.Range("D27").Value = "To test"
$ws.Range("E27").Value = "to set a bubble to advertise the users"
$ws.Range("F27").Value = "to see whether there is a error on search"
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 1

# --- Row 28: what used to be row 27's task ---
$ws.Range("C28").Value = "Low"
$ws.Range("D28").Value = "To do"
$ws.Range("E28").Value = "to improve the variable's visibility"
$ws.Range("F28").Value = "to display more variables"
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = $null

# --- Row 29: what used to be row 28's task (row was blank before) ---
$ws.Range("B29").Formula = "=ROW()"
$ws.Range("C29").Value = "Low"
$ws.Range("D29").Value = "To do"
$ws.Range("E29").Value = "to save as a list of locked variable"
$ws.Range("F29").Value = "to load quicly a list of locked elements"
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = $null

# --- Rows now "Done" drop out of the To do / To test autofilter view ---
$ws.Rows.Item(22).Hidden = $true
$ws.Rows.Item(25).Hidden = $true
$ws.Rows.Item(26).Hidden = $true

# --- Restore the cursor position left by the editing session ---
$ws.Range("E43").Select() | Out-Null
